# Weekly update: insert a new price observation row for "Poroto granado"
# at Terminal Hortofrutícola Agro Chillán, pushing all existing rows
# from 60 downward by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; this shifts rows 60-110 down to 61-111
$ws.Rows.Item(60).Insert()

# Seed the new row 60 with a copy of the row that was just pushed down to 61
# (same Mercado/Región/Categoría/Variedad/Calidad/Unidad/Origen/Clasificación)
$ws.Range("A60:R60").Value2 = $ws.Range("A61:R61").Value2

# Now overwrite the fields that differ for this new weekly observation
$ws.Range("D60").Value2 = 44942   # Fecha
$ws.Range("J60").Value2 = 60      # Volumen
$ws.Range("K60").Value2 = 40000   # Precio mínimo
$ws.Range("L60").Value2 = 40000   # Precio máximo
$ws.Range("M60").Value2 = 40000   # Precio promedio ponderado
$ws.Range("P60").Value2 = 1600    # Precio $/Kg
